$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# The source data set ("Hortaliza, Femacal de La Calera - Poroto granado")
# gained two new weekly price records. Both are inserted in their
# chronological slot among the existing rows, pushing the rows that come
# after each insertion point down by one.
#
#   1) New record dated 2021-12-24 is inserted as row 56
#      (pushing the former rows 56..115 down to 57..116)
#   2) New record dated 2021-12-23 is inserted as row 80
#      (pushing the former rows 80..116 down to 81..117)
# -------------------------------------------------------------------------

function Set-DataRow {
    param($ws, $r, $values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}

# Insert the first new row at position 56
$ws.Rows.Item(56).Insert()

$row56 = @(3, "Femacal de La Calera", "Coquimbo", "12/24/2021", 5, 100112030, "Poroto granado", "Sin especificar", "Primera", 35, 42000, 42000, 42000, "`$/saco 25 kilos", "Provincia de Limarí", 1680, 25, "Hortaliza")
Set-DataRow $ws 56 $row56
$ws.Cells.Item(56,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Insert the second new row at position 80 (post first-insert numbering)
$ws.Rows.Item(80).Insert()

$row80 = @(3, "Femacal de La Calera", "Coquimbo", "12/23/2021", 5, 100112030, "Poroto granado", "Sin especificar", "Primera", 73, 41000, 42000, 41479, "`$/malla 25 kilos", "Provincia de Limarí", 1659, 25, "Hortaliza")
Set-DataRow $ws 80 $row80
$ws.Cells.Item(80,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
